$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column F (dSF)
$updates = @{
    2  = 5
    12 = -1
    14 = -2
    19 = 6
    23 = -7
    24 = -2
    25 = -2
    29 = -4
    32 = 1
    36 = -4
    37 = 0
    42 = 5
    43 = 1
    51 = 3
    52 = -2
    53 = -2
    58 = -10
    59 = 1
    61 = 7
    65 = -1
    66 = 4
    71 = 1
    82 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
